$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values after the edit (rows 67,68,69,70,72 have their data
# rotated among each other, per the supplied diff).

# Row 67
$ws.Range("A67").Value = 111871585
$ws.Range("B67").Value = 89405
$ws.Range("D67").Value = "NT"
$ws.Range("E67").Value = 1202
$ws.Range("F67").Value = "Ullticka"
$ws.Range("G67").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H67").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P67").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q67").Value = 590630.2636057099
$ws.Range("R67").Value = 7040266.929520278

# Row 68
$ws.Range("A68").Value = 111870139
$ws.Range("B68").Value = 89845
$ws.Range("D68").Value = "VU"
$ws.Range("E68").Value = 1209
$ws.Range("F68").Value = "Rynkskinn"
$ws.Range("G68").Value = "Phlebia centrifuga"
$ws.Range("H68").Value = "P.Karst."
$ws.Range("P68").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q68").Value = 590710.4131779457
$ws.Range("R68").Value = 7040581.765558361

# Row 69
$ws.Range("A69").Value = 111870990
$ws.Range("B69").Value = 90666
$ws.Range("D69").Value = "LC"
$ws.Range("E69").Value = 4364
$ws.Range("F69").Value = "Dropptaggsvamp"
$ws.Range("G69").Value = "Hydnellum ferrugineum"
$ws.Range("H69").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P69").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q69").Value = 590569.8478412227
$ws.Range("R69").Value = 7040376.109235858

# Row 70
$ws.Range("A70").Value = 111881310
$ws.Range("B70").Value = 89425
$ws.Range("D70").Value = "NT"
$ws.Range("E70").Value = 5442
$ws.Range("F70").Value = "Tallticka"
$ws.Range("G70").Value = "Porodaedalea pini"
$ws.Range("H70").Value = "(Brot.) Murrill"
$ws.Range("P70").Value = "Valforsen, Ång"
$ws.Range("Q70").Value = 590738.9206925276
$ws.Range("R70").Value = 7040524.002523924

# Row 72
$ws.Range("A72").Value = 111870127
$ws.Range("B72").Value = 89405
$ws.Range("D72").Value = "NT"
$ws.Range("E72").Value = 1202
$ws.Range("F72").Value = "Ullticka"
$ws.Range("G72").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H72").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P72").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q72").Value = 590710.4131779457
$ws.Range("R72").Value = 7040581.765558361
